# "First set of edits after R&R"
#
# The SS sheet had a two-row header:
#   row 3 (merged C3:E3) -> "Subj. pr. of recovery"  (shared string "Forced"... no: the label text)
#   row 4 -> the real column headers (Choice / Subj. pr. of recovery / Control / p-value)
#
# The edit collapses this into a single header row: the now-redundant merged
# label row is removed, and the column header that used to read
# "Subj. pr. of recovery" is renamed to "Structure". Every row below shifts
# up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("SS")

# Remove the old merged label row (old row 3); everything below shifts up.
$ws.Rows(3).Delete()

# The old row 4 header is now row 3. Rename its "Subj. pr. of recovery"
# column header (column C) to "Structure".
$ws.Range("C3").Value = "Structure"

# That header row no longer needs the taller custom row height it had when
# it held the two-line merged caption; let it size back to the sheet default.
$ws.Rows(3).AutoFit()
